$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (columns B,C,D,E,F,I,J,K,L,M)
$data = @{
    2 = @{ "B"=1.02; "C"=1.041283512961459; "D"=1.043523996761666; "E"=1.039520552964207; "F"=1.040050791433821; "I"=1.039689566532418; "J"=1.046365555210747; "K"=1.046297620559727; "L"=1.042305506525004; "M"=1.042834238978155 }
    3 = @{ "B"=1.02; "C"=1.042611367692881; "D"=1.044557176498725; "E"=1.040663462255593; "F"=1.042007747970504; "I"=1.040095466089857; "J"=1.047337490221681; "K"=1.047141272094568; "L"=1.043257774557176; "M"=1.044598523683211 }
    4 = @{ "B"=1.019999999999999; "C"=1.043469050889066; "D"=1.045224334847763; "E"=1.041401891794451; "F"=1.043272232234967; "I"=1.04035613510211; "J"=1.047964418026637; "K"=1.04768519380586; "L"=1.043872271983033; "M"=1.045737944252976 }
    5 = @{ "B"=1.02; "C"=1.043829261698173; "D"=1.045504482618463; "E"=1.041712066810588; "F"=1.043803405673267; "I"=1.040465249799802; "J"=1.048227509637834; "K"=1.047913389285165; "L"=1.044130208137158; "M"=1.04621644559679 }
    6 = @{ "B"=1.02; "C"=1.043889721802199; "D"=1.045551501692961; "E"=1.041764131407692; "F"=1.04389256807288; "I"=1.040483543110529; "J"=1.048271656527892; "K"=1.047951676915581; "L"=1.044173493510832; "M"=1.046296758488863 }
    7 = @{ "B"=1.02; "C"=1.043473865446933; "D"=1.045228079470073; "E"=1.041406037386921; "F"=1.043279331416089; "I"=1.040357594944177; "J"=1.047967935305716; "K"=1.047688244803588; "L"=1.043875720099082; "M"=1.045744340003311 }
    8 = @{ "B"=1.02; "C"=1.041732587914412; "D"=1.043873452621659; "E"=1.039907037066103; "F"=1.040712535034915; "I"=1.039827152543522; "J"=1.046694437995569; "K"=1.04658314845425; "L"=1.042627681304074; "M"=1.043430948557376 }
    9 = @{ "B"=1.02; "C"=1.038652241919603; "D"=1.041475678783609; "E"=1.037256896222566; "F"=1.036175090221768; "I"=1.038877214315632; "J"=1.044434991111931; "K"=1.044620497051296; "L"=1.040415379287336; "M"=1.039337096938425 }
    10 = @{ "B"=1.02; "C"=1.036590202625463; "D"=1.039869661844941; "E"=1.03548397461946; "F"=1.033139512605955; "I"=1.038233535140552; "J"=1.042918035364776; "K"=1.043301491309166; "L"=1.038931417905413; "M"=1.036595362593582 }
    11 = @{ "B"=1.02; "C"=1.035695214086735; "D"=1.039172397650026; "E"=1.034714750064465; "F"=1.031822351539391; "I"=1.037952318176817; "J"=1.042258577816107; "K"=1.042727780304241; "L"=1.038286624803854; "M"=1.035405011740869 }
    12 = @{ "B"=1.02; "C"=1.03536245013538; "D"=1.038913119561226; "E"=1.034428788770395; "F"=1.031332670528145; "I"=1.037847483211902; "J"=1.042013228225196; "K"=1.042514286731846; "L"=1.038046779540654; "M"=1.034962371426021 }
    13 = @{ "B"=1.02; "C"=1.035433843986005; "D"=1.038968748491191; "E"=1.034490139262548; "F"=1.031437728454886; "I"=1.037869987841229; "J"=1.042065874642313; "K"=1.042560099640275; "L"=1.03809824265887; "M"=1.03505734175349 }
    14 = @{ "B"=1.02; "C"=1.035667714398852; "D"=1.039150971449581; "E"=1.034691117285293; "F"=1.031781883203502; "I"=1.037943660223457; "J"=1.042238305291426; "K"=1.042710140881154; "L"=1.038266806096487; "M"=1.035368433063775 }
    15 = @{ "B"=1.02; "C"=1.035811766323142; "D"=1.039263207348343; "E"=1.034814914886769; "F"=1.03199387095742; "I"=1.037989001971669; "J"=1.042344492619781; "K"=1.042802534172369; "L"=1.038370618317684; "M"=1.035560041143289 }
    16 = @{ "B"=1.02; "C"=1.036649554884635; "D"=1.039915897569308; "E"=1.035534992608689; "F"=1.033226869097917; "I"=1.038252145671138; "J"=1.042961745961935; "K"=1.043339511968103; "L"=1.038974163276145; "M"=1.036674294283109 }
    17 = @{ "B"=1.02; "C"=1.037174506422346; "D"=1.040324814245112; "E"=1.035986262403655; "F"=1.033999551324973; "I"=1.038416537485011; "J"=1.043348230684766; "K"=1.043675651464849; "L"=1.039352150597754; "M"=1.037372378100631 }
    18 = @{ "B"=1.02; "C"=1.037480498552709; "D"=1.040563150464379; "E"=1.03624933247421; "F"=1.034449981021944; "I"=1.038512183469583; "J"=1.043573409696201; "K"=1.04387146806312; "L"=1.039572409620497; "M"=1.03777925460121 }
    19 = @{ "B"=1.02; "C"=1.037584799734647; "D"=1.04064438693444; "E"=1.03633900765288; "F"=1.034603521870611; "I"=1.038544755492169; "J"=1.043650147534905; "K"=1.043938194574525; "L"=1.039647476027716; "M"=1.037917937894384 }
    20 = @{ "B"=1.02; "C"=1.037118205170542; "D"=1.040280959803586; "E"=1.03593786077511; "F"=1.033916677145479; "I"=1.03839892474661; "J"=1.043306790538944; "K"=1.043639612556753; "L"=1.039311618373345; "M"=1.037297511873991 }
    21 = @{ "B"=1.02; "C"=1.035598854437976; "D"=1.039097319189624; "E"=1.034631940858973; "F"=1.031680550186303; "I"=1.037921976014176; "J"=1.042187539791545; "K"=1.042665968332122; "L"=1.03821717780776; "M"=1.035276838109974 }
    22 = @{ "B"=1.02; "C"=1.034641691856669; "D"=1.0383514765705; "E"=1.033809482690573; "F"=1.030272118944647; "I"=1.037619908275768; "J"=1.041481518360322; "K"=1.042051530332815; "L"=1.037527087317497; "M"=1.034003511470125 }
    23 = @{ "B"=1.02; "C"=1.035149283378646; "D"=1.038747019212882; "E"=1.034245615633045; "F"=1.031018996920252; "I"=1.0377802487767; "J"=1.041856014169128; "K"=1.042377472396081; "L"=1.037893106195083; "M"=1.034678801426717 }
    24 = @{ "B"=1.02; "C"=1.037143645910016; "D"=1.040300776287466; "E"=1.03595973184435; "F"=1.033954125230034; "I"=1.038406883929307; "J"=1.043325516332057; "K"=1.043655897753532; "L"=1.03932993380344; "M"=1.037331341636849 }
    25 = @{ "B"=1.02; "C"=1.039450048590259; "D"=1.042096862610072; "E"=1.037943084584007; "F"=1.037349931550898; "I"=1.03912461606679; "J"=1.045020967466598; "K"=1.045129733255931; "L"=1.040988893958362; "M"=1.040397597793817 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($c in $rowVals.Keys) {
        $ws.Range("$c$r").Value = $rowVals[$c]
    }
}

Write-Host "Updated $($data.Count) rows"
